$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.803.22"
$ws.Range("E2").Value = "  -4.36%  "
$ws.Range("D3").Value = "3.498.50"
$ws.Range("E3").Value = "  -5.24%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "579.32"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").Value = "174.11"
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "3.491.33"
$ws.Range("E8").Value = "  -5.23%  "
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "0.188"
$ws.Range("E10").Value = "  -7.56%  "
$ws.Range("E11").Value = "  +9.42%  "
$ws.Range("D12").Value = "0.600"
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("D13").Value = "47.13"
$ws.Range("E13").Value = "  -5.76%  "
$ws.Range("D14").Value = "0.0000276"
$ws.Range("E14").Value = "  -3.88%  "
$ws.Range("D15").Value = "671.10"
$ws.Range("E15").Value = "  -2.25%  "
$ws.Range("D16").Value = "4.062.80"
$ws.Range("E16").Value = "  -5.22%  "
$ws.Range("D17").Value = "8.82"
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("D18").Value = "3.502.80"
$ws.Range("E18").Value = "  -5.16%  "
$ws.Range("D19").Value = "68.760.61"
$ws.Range("E19").Value = "  -4.57%  "
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").Value = "17.50"
$ws.Range("E21").Value = "  -4.49%  "
$ws.Range("D22").Value = "11.16"
$ws.Range("E22").Value = "  -4.11%  "
$ws.Range("D23").Value = "0.902"
$ws.Range("E23").Value = "  -4.37%  "
$ws.Range("D24").Value = "16.26"
$ws.Range("E24").Value = "  -8.68%  "
$ws.Range("D25").Value = "98.05"
$ws.Range("E25").Value = "  -5.57%  "
$ws.Range("D26").Value = "3.87"
$ws.Range("E26").Value = "  -4.26%  "
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "2.64"
$ws.Range("E29").Value = "  -7.51%  "
$ws.Range("D30").Value = "9.41"
$ws.Range("E30").Value = "  -8.15%  "
$ws.Range("D31").Value = "32.85"
$ws.Range("E31").Value = "  -7.80%  "
$ws.Range("D32").Value = "8.70"
$ws.Range("E32").Value = "  -5.53%  "
$ws.Range("E33").Value = "  -8.29%  "
$ws.Range("D34").Value = "7.28"
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("E35").Value = "  -6.08%  "
$ws.Range("D36").Value = "578.30"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("B37").Value = "Cosmos"
$ws.Range("C37").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D37").Value = "10.92"
$ws.Range("E37").Value = "  -3.78%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "3.58"
$ws.Range("E38").Value = "  -14.53%  "
$ws.Range("E39").Value = "  -4.40%  "
$ws.Range("E40").Value = "  -5.28%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.137"
$ws.Range("E42").Value = "  -5.46%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "0.336"
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0437"
$ws.Range("E44").Value = "  -5.50%  "
$ws.Range("D45").Value = "3.411.86"
$ws.Range("E45").Value = "  -9.30%  "
$ws.Range("D46").Value = "33.29"
$ws.Range("E46").Value = "  -6.52%  "
$ws.Range("D47").Value = "0.0₃0701"
$ws.Range("E47").Value = "  -9.71%  "
$ws.Range("D48").Value = "2.87"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("E49").Value = "  -7.62%  "
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").Value = "132.14"
$ws.Range("E51").Value = "  -1.49%  "
